# Updating results with correct MI SEs
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# theta_se row (row 4)
$ws.Range("B4").Value = "(1.11)"
$ws.Range("C4").Value = "(0.51)"

# lambda_se row (row 6)
$ws.Range("B6").Value = "(0.97)"
$ws.Range("C6").Value = "(0.42)"
